$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 3rd data row
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-07 07:04:52"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 3rd data row
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-07 07:04:47"
$wsZhCn.Range("K4").Value = "2016-09-07 07:05:17"

# de-de sheet: the same "Latest HO Xliff Generate Date" value is shared with Overview's G4,
# so it also needs to be refreshed to stay in sync; Correspond Handback DateTime is updated too
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-07 07:04:52"
$wsDeDe.Range("K4").Value = "2016-09-07 07:05:26"
